# Chiffres COVID-19 Valais.xlsx - daily data update
# Adds/updates the daily COVID figures for Valais and moves the frozen-pane
# view to show the most recently edited rows (matching the author's
# workflow of scrolling back up to row 3 then re-selecting the new last
# entry row in the bottom-right pane).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data corrections on already-entered rows -------------------------

# Row 254: new positive cases revised 543 -> 544
$ws.Range("C254").Value = 544

# Rows 517-521: "Patients aux SI" (G) corrected down by one patient each day
$ws.Range("G517").Value = 1
$ws.Range("G518").Value = 1
$ws.Range("G519").Value = 2
$ws.Range("G520").Value = 2
$ws.Range("G521").Value = 2

# Row 522: new positive cases revised 17 -> 16
$ws.Range("C522").Value = 16

# --- New daily entries --------------------------------------------------

# Row 524 (2021-08-02): new positive cases
$ws.Range("C524").Value = 33

# Row 525 (2021-08-03)
$ws.Range("C525").Value = 31
$ws.Range("E525").Value = 1
$ws.Range("F525").Value = 1
$ws.Range("G525").Value = 2
$ws.Range("L525").Value = 0
$ws.Range("M525").Value = 0

# Row 526 (2021-08-04)
$ws.Range("C526").Value = 21
$ws.Range("E526").Value = 1
$ws.Range("F526").Value = 1
$ws.Range("G526").Value = 2
$ws.Range("L526").Value = 0
$ws.Range("M526").Value = 0

# Row 527 (2021-08-05)
$ws.Range("C527").Value = 2
$ws.Range("E527").Value = 1
$ws.Range("F527").Value = 1
$ws.Range("G527").Value = 2
$ws.Range("L527").Value = 0
$ws.Range("M527").Value = 0

# --- View state: scroll frozen pane back to top, select new last row ---

$ws.Range("B3").Select()
$ws.Range("O521").Select()
